# Fixed traj read real robot data
# Update B:F columns for rows 1-6 on the active sheet with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = @(1.732704007046913, 0.1945306715051764, -0.672838158291254, 0.7034277224914169, 1.570796292848413)
    2 = @(1.778290465355663, 0.1934420031981677, -0.6758680428538398, 0.7014865303281086, 1.570796289207822)
    3 = @(1.982548021645145, 0.1885640474137254, -0.6894439344574702, 0.6927887024582571, 1.570796272895559)
    4 = @(2.269919010210249, 0.1817012267601512, -0.708543924987362, 0.680551684456407, 1.570796249945755)
    5 = @(2.474176566499732, 0.1768232709757089, -0.7221198165909923, 0.6718538565865555, 1.570796233633492)
    6 = @(2.519763024808483, 0.1757346026687003, -0.7251497011535781, 0.6699126644232472, 1.570796229992901)
}

foreach ($row in 1..6) {
    $rowValues = $values[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = 2 + $i  # Column B = 2, C = 3, D = 4, E = 5, F = 6
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}
